$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header styling (A1 already has the bold/border header style) to the
#     new header cells G1:N1 before writing their text, so the style survives. ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:N1").PasteSpecial(-4122) | Out-Null

# --- Row 1: headers (A1:N1) ---
$headers = @{
    1  = "Tipo"
    2  = "Nome"
    3  = "Cargo"
    4  = "Data"
    5  = "conquistas"
    6  = "desafios"
    7  = "autoavaliacao"
    8  = "objetivos"
    9  = "gestor"
    10 = "colaborador"
    11 = "pontos_fortes"
    12 = "areas_desenvolvimento"
    13 = "avaliacao"
    14 = "plano_desenvolvimento"
}
foreach ($col in $headers.Keys) {
    $ws.Cells.Item(1, $col).Value = $headers[$col]
}

# --- Data rows 2-6. Use column-letter keyed hashtables; $null means "leave blank". ---
$rows = @(
    @{ A="auto"; B="iago"; C="desenvolvedor full-stack"; D="06/06/2025 14:56"; E="aaaaaaaa"; F="aaaaaaaaaa"; G="4"; H="aaaaaaaaaaaaaa" },
    @{ A="auto"; D="06/06/2025 18:00"; I="iago"; J="Carlin"; K="aaaaaaaaaaaa"; L="aaaaaaaaaaaaaa"; M="4"; N="aaaaaaaaaaaaaaaaaa" },
    @{ A="auto"; B="ana"; C="desenvolvedor full-stack"; D="06/06/2025 22:18"; E="aaaaa"; F="aaaaaaaaaaa"; G="3"; H="aaaaaaaaaaaaaaaa" },
    @{ A="auto"; B="gabi"; C="design"; D="06/06/2025 22:21"; E="aaaa"; F="aaaaaaaa"; G="5"; H="iago" },
    @{ A="auto"; D="06/06/2025 22:22"; I="gabi"; J="iago"; K="teste"; L="teste"; M="5"; N="teste 001" }
)

$colIndex = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14 }
# Columns whose value is a plain-digit string and must stay *text*, not be
# coerced into a number by Excel's input parser.
$numericLooking = @{ "2:G"=1; "3:M"=1; "4:G"=1; "5:G"=1; "6:M"=1 }

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    foreach ($col in $colIndex.Keys) {
        $cIdx = $colIndex[$col]
        $val = $rowData[$col]
        if ($null -eq $val) { continue }
        $key = "$($rowNum):$col"
        if ($numericLooking.ContainsKey($key)) {
            # Force text storage (quote-prefix), then strip the resulting
            # "quoted text" style back to the sheet default so no extra
            # per-cell style is introduced.
            $ws.Cells.Item($rowNum, $cIdx).Formula = "'" + $val
            $ws.Cells.Item($rowNum, $cIdx).Style = "Normal"
        } else {
            $ws.Cells.Item($rowNum, $cIdx).Value = $val
        }
    }
}
